$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.877.21"
$ws.Range("E2").Value = "  +1.09%  "
# Row 3
$ws.Range("D3").Value = "1.639.14"
$ws.Range("E3").Value = "  -0.35%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.39%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.43%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.524"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.47%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.38%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.39"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.29%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.261"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.83%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0611"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.18%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0891"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.24%  "
# Row 12
$ws.Range("D12").Value = "1.869.53"
$ws.Range("E12").Value = "  -0.52%  "
# Row 13
$ws.Range("D13").Value = "1.617.13"
$ws.Range("E13").Value = "  -1.81%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.44%  "
# Row 15
$ws.Range("E15").Value = "  -3.89%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.16%  "
# Row 17
$ws.Range("D17").Value = "27.886.33"
$ws.Range("E17").Value = "  +1.27%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "231.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.90%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.12%  "
# Row 20
$ws.Range("D20").Value = "0.0₃0723"
$ws.Range("E20").Value = "  +0.06%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.37%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.11%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.12%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.42%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.34%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.82%  "
# Row 27
$ws.Range("E27").Value = "  -0.68%  "
# Row 28
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.34%  "
# Row 29
$ws.Range("B29").Value = "BinanceUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.34%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.08%  "
# Row 31
$ws.Range("E31").Value = "  -0.91%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.31"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.19%  "
# Row 33
$ws.Range("D33").Value = "1.473.54"
$ws.Range("E33").Value = "  +3.29%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.30%  "
# Row 35
$ws.Range("E35").Value = "  -2.73%  "
# Row 36
$ws.Range("E36").Value = "  -1.02%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.568"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.61%  "
# Row 38
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.928"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +13.59%  "
# Row 39
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.878"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.52%  "
# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0167"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.20%  "
# Row 41
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "68.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.68%  "
# Row 42
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.38%  "
# Row 43
$ws.Range("E43").Value = "  -2.41%  "
# Row 44
$ws.Range("E44").Value = "  -0.45%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.80%  "
# Row 46
$ws.Range("E46").Value = "  -0.43%  "
# Row 47
$ws.Range("D47").Value = "1.780.43"
$ws.Range("E47").Value = "  -0.45%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.37%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "87.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.02%  "
# Row 50
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0992"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.60%  "
# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.82%  "
